# This script implements the "added harvard case classification" edit:
# the workbook gains a new doctor-average column. The previously computed
# "average_doctor" figures (column BP) become the new "average_doctor_old"
# baseline (column BQ), the BP/BQ headers swap accordingly, column BP is
# refreshed with the new averages, and every per-app variance/std-dev/avg
# stat on rows 4-13 is recomputed against the updated case set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header swap: BP was "average_doctor" / BQ was "average_doctor_old"
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"


# Row 4: stats_for_precision
$ws.Range("E4").Value = 0.377
$ws.Range("F4").Value = 0.08400000000000001
$ws.Range("G4").Value = 0.289
$ws.Range("N4").Value = 0.396
$ws.Range("O4").Value = 0.064
$ws.Range("P4").Value = 0.253
$ws.Range("W4").Value = 0.212
$ws.Range("X4").Value = 0.1
$ws.Range("Y4").Value = 0.317
$ws.Range("AI4").Value = 0.193
$ws.Range("AJ4").Value = 0.064
$ws.Range("AK4").Value = 0.253
$ws.Range("AU4").Value = 0.145
$ws.Range("AV4").Value = 0.026
$ws.Range("AW4").Value = 0.161
$ws.Range("BA4").Value = 1.969
$ws.Range("BB4").Value = 0.17
$ws.Range("BC4").Value = 0.412
$ws.Range("BG4").Value = 0.735
$ws.Range("BH4").Value = 0.141
$ws.Range("BI4").Value = 0.375
$ws.Range("BM4").Value = 0.6879999999999999
$ws.Range("BN4").Value = 0.091
$ws.Range("BO4").Value = 0.302
$ws.Range("BP4").Value = 0.656
$ws.Range("BQ4").Value = 0.662

# Row 5: stats_for_recall
$ws.Range("E5").Value = 0.487
$ws.Range("F5").Value = 0.102
$ws.Range("G5").Value = 0.32
$ws.Range("N5").Value = 0.747
$ws.Range("O5").Value = 0.08500000000000001
$ws.Range("P5").Value = 0.292
$ws.Range("W5").Value = 0.221
$ws.Range("X5").Value = 0.112
$ws.Range("Y5").Value = 0.335
$ws.Range("AI5").Value = 0.231
$ws.Range("AJ5").Value = 0.092
$ws.Range("AK5").Value = 0.304
$ws.Range("AU5").Value = 0.298
$ws.Range("AV5").Value = 0.097
$ws.Range("AW5").Value = 0.312
$ws.Range("BA5").Value = 1.368
$ws.Range("BB5").Value = 0.08400000000000001
$ws.Range("BC5").Value = 0.29
$ws.Range("BG5").Value = 0.411
$ws.Range("BH5").Value = 0.051
$ws.Range("BM5").Value = 0.572
$ws.Range("BP5").Value = 0.456
$ws.Range("BQ5").Value = 0.455

# Row 6: stats_for_f1-score
$ws.Range("E6").Value = 0.425
$ws.Range("N6").Value = 0.518
$ws.Range("W6").Value = 0.216
$ws.Range("AI6").Value = 0.21
$ws.Range("AU6").Value = 0.195
$ws.Range("BA6").Value = 1.604
$ws.Range("BG6").Value = 0.527
$ws.Range("BM6").Value = 0.625
$ws.Range("BP6").Value = 0.535
$ws.Range("BQ6").Value = 0.536

# Row 7: stats_for_f2-score
$ws.Range("E7").Value = 0.46
$ws.Range("N7").Value = 0.635
$ws.Range("W7").Value = 0.219
$ws.Range("AI7").Value = 0.222
$ws.Range("AU7").Value = 0.246
$ws.Range("BA7").Value = 1.452
$ws.Range("BG7").Value = 0.451
$ws.Range("BM7").Value = 0.592
$ws.Range("BQ7").Value = 0.484

# Row 8: stats_for_NDCG
$ws.Range("E8").Value = 0.519
$ws.Range("F8").Value = 0.13
$ws.Range("G8").Value = 0.361
$ws.Range("N8").Value = 0.748
$ws.Range("O8").Value = 0.07099999999999999
$ws.Range("P8").Value = 0.266
$ws.Range("W8").Value = 0.213
$ws.Range("X8").Value = 0.104
$ws.Range("Y8").Value = 0.323
$ws.Range("AI8").Value = 0.211
$ws.Range("AJ8").Value = 0.091
$ws.Range("AK8").Value = 0.302
$ws.Range("AU8").Value = 0.236
$ws.Range("AV8").Value = 0.074
$ws.Range("AW8").Value = 0.273
$ws.Range("BA8").Value = 1.731
$ws.Range("BB8").Value = 0.135
$ws.Range("BC8").Value = 0.368
$ws.Range("BG8").Value = 0.571
$ws.Range("BH8").Value = 0.104
$ws.Range("BI8").Value = 0.323
$ws.Range("BM8").Value = 0.704
$ws.Range("BN8").Value = 0.07099999999999999
$ws.Range("BO8").Value = 0.266
$ws.Range("BP8").Value = 0.577
$ws.Range("BQ8").Value = 0.587

# Row 9: stats_for_M1
$ws.Range("E9").Value = 0.451
$ws.Range("F9").Value = 0.248
$ws.Range("G9").Value = 0.498
$ws.Range("N9").Value = 0.627
$ws.Range("O9").Value = 0.234
$ws.Range("P9").Value = 0.483
$ws.Range("W9").Value = 0.118
$ws.Range("X9").Value = 0.104
$ws.Range("Y9").Value = 0.322
$ws.Range("AI9").Value = 0.118
$ws.Range("AJ9").Value = 0.104
$ws.Range("AK9").Value = 0.322
$ws.Range("BA9").Value = 1.686
$ws.Range("BB9").Value = 0.245
$ws.Range("BC9").Value = 0.495
$ws.Range("BG9").Value = 0.608
$ws.Range("BH9").Value = 0.238
$ws.Range("BI9").Value = 0.488
$ws.Range("BM9").Value = 0.647
$ws.Range("BN9").Value = 0.228
$ws.Range("BO9").Value = 0.478
$ws.Range("BP9").Value = 0.5620000000000001
$ws.Range("BQ9").Value = 0.5610000000000001

# Row 10: stats_for_M3
$ws.Range("E10").Value = 0.569
$ws.Range("F10").Value = 0.245
$ws.Range("G10").Value = 0.495
$ws.Range("N10").Value = 0.824
$ws.Range("O10").Value = 0.145
$ws.Range("P10").Value = 0.381
$ws.Range("W10").Value = 0.255
$ws.Range("X10").Value = 0.19
$ws.Range("Y10").Value = 0.436
$ws.Range("AI10").Value = 0.235
$ws.Range("AJ10").Value = 0.18
$ws.Range("AK10").Value = 0.424
$ws.Range("AU10").Value = 0.216
$ws.Range("AV10").Value = 0.169
$ws.Range("AW10").Value = 0.411
$ws.Range("BA10").Value = 2
$ws.Range("BB10").Value = 0.25
$ws.Range("BC10").Value = 0.5
$ws.Range("BG10").Value = 0.647
$ws.Range("BH10").Value = 0.228
$ws.Range("BI10").Value = 0.478
$ws.Range("BM10").Value = 0.863
$ws.Range("BN10").Value = 0.118
$ws.Range("BO10").Value = 0.344
$ws.Range("BP10").Value = 0.667
$ws.Range("BQ10").Value = 0.6929999999999999

# Row 11: stats_for_M5
$ws.Range("E11").Value = 0.588
$ws.Range("F11").Value = 0.242
$ws.Range("G11").Value = 0.492
$ws.Range("N11").Value = 0.863
$ws.Range("O11").Value = 0.118
$ws.Range("P11").Value = 0.344
$ws.Range("W11").Value = 0.255
$ws.Range("X11").Value = 0.19
$ws.Range("Y11").Value = 0.436
$ws.Range("AI11").Value = 0.235
$ws.Range("AJ11").Value = 0.18
$ws.Range("AK11").Value = 0.424
$ws.Range("AU11").Value = 0.333
$ws.Range("AV11").Value = 0.222
$ws.Range("AW11").Value = 0.471
$ws.Range("BA11").Value = 2
$ws.Range("BB11").Value = 0.25
$ws.Range("BC11").Value = 0.5
$ws.Range("BG11").Value = 0.647
$ws.Range("BH11").Value = 0.228
$ws.Range("BI11").Value = 0.478
$ws.Range("BM11").Value = 0.863
$ws.Range("BN11").Value = 0.118
$ws.Range("BO11").Value = 0.344
$ws.Range("BP11").Value = 0.667
$ws.Range("BQ11").Value = 0.6929999999999999

# Row 12: stats_for_position
$ws.Range("E12").Value = 1.4
$ws.Range("F12").Value = 0.64
$ws.Range("G12").Value = 0.8
$ws.Range("N12").Value = 1.652
$ws.Range("O12").Value = 1.618
$ws.Range("P12").Value = 1.272
$ws.Range("W12").Value = 1.846
$ws.Range("X12").Value = 0.746
$ws.Range("Y12").Value = 0.863
$ws.Range("AI12").Value = 1.917
$ws.Range("AJ12").Value = 0.91
$ws.Range("AK12").Value = 0.954
$ws.Range("AU12").Value = 2.647
$ws.Range("AV12").Value = 1.758
$ws.Range("AW12").Value = 1.326
$ws.Range("BB12").Value = 0.32
$ws.Range("BC12").Value = 0.5659999999999999
$ws.Range("BG12").Value = 1.061
$ws.Range("BH12").Value = 0.057
$ws.Range("BI12").Value = 0.239
$ws.Range("BM12").Value = 1.341
$ws.Range("BN12").Value = 0.407
$ws.Range("BO12").Value = 0.638
$ws.Range("BP12").Value = 1.201
$ws.Range("BQ12").Value = 1.263

# Row 13: stats_for_length (x of gs)
$ws.Range("E13").Value = 1.714
$ws.Range("F13").Value = 0.912
$ws.Range("G13").Value = 0.955
$ws.Range("N13").Value = 2.341
$ws.Range("O13").Value = 1.144
$ws.Range("P13").Value = 1.07
$ws.Range("W13").Value = 1.09
$ws.Range("X13").Value = 0.186
$ws.Range("Y13").Value = 0.431
$ws.Range("AI13").Value = 1.383
$ws.Range("AJ13").Value = 0.401
$ws.Range("AK13").Value = 0.633
$ws.Range("AU13").Value = 2.482
$ws.Range("AV13").Value = 1.307
$ws.Range("AW13").Value = 1.143
$ws.Range("BA13").Value = 2.5
$ws.Range("BB13").Value = 0.31
$ws.Range("BC13").Value = 0.5570000000000001
$ws.Range("BG13").Value = 0.619
$ws.Range("BH13").Value = 0.08500000000000001
$ws.Range("BI13").Value = 0.292
$ws.Range("BM13").Value = 0.989
$ws.Range("BN13").Value = 0.361
$ws.Range("BO13").Value = 0.601
$ws.Range("BP13").Value = 0.833
$ws.Range("BQ13").Value = 0.782
